$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells that change to Text format so that
# numeric-looking values (e.g. "303.80", "71.40") are preserved exactly
# as text instead of being normalized as numbers by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23", "D25", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D42", "D43", "D46", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) column values
$ws.Range("D2").Value = "43.150.62"
$ws.Range("D3").Value = "2.372.02"
$ws.Range("D5").Value = "303.80"
$ws.Range("D6").Value = "95.63"
$ws.Range("D9").Value = "0.482"
$ws.Range("D10").Value = "34.44"
$ws.Range("D12").Value = "0.0788"
$ws.Range("D13").Value = "18.21"
$ws.Range("D14").Value = "6.78"
$ws.Range("D15").Value = "2.739.05"
$ws.Range("D16").Value = "2.362.87"
$ws.Range("D17").Value = "0.800"
$ws.Range("D18").Value = "43.167.45"
$ws.Range("D19").Value = "11.98"
$ws.Range("D22").Value = "67.94"
$ws.Range("D23").Value = "235.33"
$ws.Range("D25").Value = "2.45"
$ws.Range("D27").Value = "24.51"
$ws.Range("D28").Value = "2.37"
$ws.Range("D29").Value = "9.35"
$ws.Range("D30").Value = "32.14"
$ws.Range("D32").Value = "5.03"
$ws.Range("D33").Value = "17.70"
$ws.Range("D34").Value = "0.110"
$ws.Range("D35").Value = "0.0734"
$ws.Range("D36").Value = "128.87"
$ws.Range("D42").Value = "21.13"
$ws.Range("D43").Value = "1.930.43"
$ws.Range("D46").Value = "9.24"
$ws.Range("D48").Value = "2.599.35"
$ws.Range("D50").Value = "71.40"
$ws.Range("D51").Value = "51.50"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +9.48%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("E46").Value = "  -8.69%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  -2.41%  "

